$wb = $excel.ActiveWorkbook

$sheet2021 = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1) Duplicate the "总计" sheet. The duplicate keeps all of its original
#    formatting/styles and will become the new, final "总计" sheet (we
#    add the extra "2022-Q1" summary row to it below). The original
#    sheet object is repurposed below into the new "2022-Q1" sheet.
# ---------------------------------------------------------------------
$totalSheet.Copy($null, $totalSheet)

$q1Sheet = $wb.Worksheets.Item("总计")
$newTotalSheet = $wb.Worksheets.Item("总计 (2)")

# Rename the original sheet out of the way first, then rename the
# duplicate back to "总计" (renaming straight to a name that is still in
# use elsewhere fails silently).
$q1Sheet.Name = "2022-Q1"
$newTotalSheet.Name = "总计"

# ---------------------------------------------------------------------
# 2) Build the "2022-Q1" sheet using the same column layout/styling as
#    "2021-Q4" (same headers, bold/bordered header row, bordered first
#    data row) by copying that sheet's used range across, then updating
#    the data that differs.
# ---------------------------------------------------------------------
$q1Sheet.Cells.Clear()
$sheet2021.Range("A1:H2").Copy($q1Sheet.Range("A1:H2"))

# Fund code/name/size are unchanged from 2021-Q4, only the position
# metrics + rank differ for the new quarter. Keep them stored as text
# (matching how 2021-Q4 stores these numeric-looking values) by forcing
# a text number format before assigning.
$q1Sheet.Range("D2:G2").NumberFormat = "@"
$q1Sheet.Range("D2").Value = "0.14"
$q1Sheet.Range("E2").Value = "83.76"
$q1Sheet.Range("F2").Value = "4.56"
$q1Sheet.Range("G2").Value = "0.0064"
$q1Sheet.Range("H2").Value = 7

# ---------------------------------------------------------------------
# 3) Insert the new "2022-Q1" row at the top of the "总计" sheet's data,
#    pushing the existing "2021-Q4" row down to row 3. Use an in-sheet
#    copy (rather than Rows.Insert) so the original cell styling carries
#    over faithfully to both rows.
# ---------------------------------------------------------------------
$newTotalSheet.Range("A2:D2").Copy($newTotalSheet.Range("A3:D3"))
$newTotalSheet.Range("A3").Value = 1

$newTotalSheet.Range("A2").Value = 0
$newTotalSheet.Range("B2").Value = "2022-Q1"
$newTotalSheet.Range("C2").Value = 1
$newTotalSheet.Range("D2").Value = 0.01
